$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $style = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $style
}

Set-TextValue ($ws.Range('D2')) '26.669.34'
Set-TextValue ($ws.Range('E2')) '  -0.43%  '
Set-TextValue ($ws.Range('D3')) '1.630.69'
Set-TextValue ($ws.Range('E3')) '  -1.08%  '
Set-TextValue ($ws.Range('E4')) '  -0.08%  '
Set-TextValue ($ws.Range('D5')) '217.69'
Set-TextValue ($ws.Range('E5')) '  +0.58%  '
Set-TextValue ($ws.Range('E6')) '  -1.72%  '
Set-TextValue ($ws.Range('E7')) '  -0.03%  '
Set-TextValue ($ws.Range('E8')) '  -1.50%  '
Set-TextValue ($ws.Range('D9')) '0.0619'
Set-TextValue ($ws.Range('E9')) '  -1.19%  '
Set-TextValue ($ws.Range('D10')) '18.94'
Set-TextValue ($ws.Range('E10')) '  -1.68%  '
Set-TextValue ($ws.Range('E11')) '  -0.18%  '
Set-TextValue ($ws.Range('D12')) '1.860.79'
Set-TextValue ($ws.Range('E12')) '  -0.98%  '
Set-TextValue ($ws.Range('D13')) '1.637.82'
Set-TextValue ($ws.Range('E13')) '  -0.83%  '
Set-TextValue ($ws.Range('E14')) '  -2.38%  '
Set-TextValue ($ws.Range('E15')) '  -2.27%  '
Set-TextValue ($ws.Range('D16')) '63.90'
Set-TextValue ($ws.Range('E16')) '  -2.37%  '
Set-TextValue ($ws.Range('D17')) '26.675.36'
Set-TextValue ($ws.Range('E17')) '  -0.44%  '
Set-TextValue ($ws.Range('E18')) '  -3.33%  '
Set-TextValue ($ws.Range('B19')) 'BitcoinCash'
Set-TextValue ($ws.Range('C19')) 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue ($ws.Range('D19')) '211.20'
Set-TextValue ($ws.Range('E19')) '  -3.22%  '
Set-TextValue ($ws.Range('B20')) 'Dai'
Set-TextValue ($ws.Range('C20')) 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue ($ws.Range('D20')) '1.01'
Set-TextValue ($ws.Range('E20')) '  -0.06%  '
Set-TextValue ($ws.Range('E21')) '  -1.90%  '
Set-TextValue ($ws.Range('D22')) '2.33'
Set-TextValue ($ws.Range('E22')) '  -7.41%  '
Set-TextValue ($ws.Range('E23')) '  -2.78%  '
Set-TextValue ($ws.Range('E24')) '  -3.55%  '
Set-TextValue ($ws.Range('D25')) '146.46'
Set-TextValue ($ws.Range('E25')) '  +0.35%  '
Set-TextValue ($ws.Range('E26')) '  +0.08%  '
Set-TextValue ($ws.Range('E27')) '  -2.80%  '
Set-TextValue ($ws.Range('D28')) '6.99'
Set-TextValue ($ws.Range('E28')) '  -2.85%  '
Set-TextValue ($ws.Range('E29')) '  -2.13%  '
Set-TextValue ($ws.Range('D30')) '0.0501'
Set-TextValue ($ws.Range('E30')) '  -3.49%  '
Set-TextValue ($ws.Range('D31')) '1.18'
Set-TextValue ($ws.Range('E31')) '  +0.43%  '
Set-TextValue ($ws.Range('D32')) '3.36'
Set-TextValue ($ws.Range('E32')) '  +0.06%  '
Set-TextValue ($ws.Range('E33')) '  -2.84%  '
Set-TextValue ($ws.Range('D34')) '1.257.76'
Set-TextValue ($ws.Range('E34')) '  -1.97%  '
Set-TextValue ($ws.Range('E35')) '  -0.05%  '
Set-TextValue ($ws.Range('E36')) '  -2.81%  '
Set-TextValue ($ws.Range('E37')) '  -3.55%  '
Set-TextValue ($ws.Range('D38')) '0.521'
Set-TextValue ($ws.Range('E38')) '  -3.77%  '
Set-TextValue ($ws.Range('E39')) '  -0.09%  '
Set-TextValue ($ws.Range('D40')) '0.801'
Set-TextValue ($ws.Range('E40')) '  -3.88%  '
Set-TextValue ($ws.Range('D41')) '0.796'
Set-TextValue ($ws.Range('E41')) '  -2.70%  '
Set-TextValue ($ws.Range('D42')) '2.16'
Set-TextValue ($ws.Range('E42')) '  -4.13%  '
Set-TextValue ($ws.Range('D43')) '1.771.85'
Set-TextValue ($ws.Range('E43')) '  -1.06%  '
Set-TextValue ($ws.Range('D44')) '5.24'
Set-TextValue ($ws.Range('E44')) '  -3.88%  '
Set-TextValue ($ws.Range('D45')) '90.82'
Set-TextValue ($ws.Range('E45')) '  -1.37%  '
Set-TextValue ($ws.Range('D46')) '59.72'
Set-TextValue ($ws.Range('E46')) '  -0.23%  '
Set-TextValue ($ws.Range('E47')) '  -2.75%  '
Set-TextValue ($ws.Range('D48')) '0.0515'
Set-TextValue ($ws.Range('E48')) '  -0.12%  '
Set-TextValue ($ws.Range('E49')) '  +0.12%  '
Set-TextValue ($ws.Range('B50')) 'Algorand'
Set-TextValue ($ws.Range('C50')) 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue ($ws.Range('D50')) '0.0953'
Set-TextValue ($ws.Range('E50')) '  -2.68%  '
Set-TextValue ($ws.Range('B51')) 'Mantle'
Set-TextValue ($ws.Range('C51')) 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue ($ws.Range('D51')) '0.405'
Set-TextValue ($ws.Range('E51')) '  -0.89%  '
